$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q4" right before the current "2022-Q3"
#    sheet, and fill it with the new fund-holding data. Cloning the
#    existing "2022-Q3" sheet (rather than Worksheets.Add()) carries over
#    its sheet-level setup (outline/page-setup props, margins, column
#    widths, header styling) so the new sheet looks consistent with its
#    siblings.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3Template = $wb.Worksheets.Item(2)
$q3Template.Copy($q3Template)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# The template (2022-Q3) had 3 fund rows; 2022-Q4 only has 2, so drop the
# extra cloned row.
$q4Sheet.Rows.Item(4).Delete()

# Row 2: fund 009394
$q4Sheet.Cells.Item(2, 1).Value = 0
$q4Sheet.Cells.Item(2, 2).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 2).Value = "009394"
$q4Sheet.Cells.Item(2, 3).Value = "银华同力精选混合"
$q4Sheet.Cells.Item(2, 4).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 4).Value = "18.05"
$q4Sheet.Cells.Item(2, 5).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 5).Value = "94.33"
$q4Sheet.Cells.Item(2, 6).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 6).Value = "5.70"
$q4Sheet.Cells.Item(2, 7).NumberFormat = "@"
$q4Sheet.Cells.Item(2, 7).Value = "1.0288"
$q4Sheet.Cells.Item(2, 8).Value = 8

# Row 3: fund 180020
$q4Sheet.Cells.Item(3, 1).Value = 1
$q4Sheet.Cells.Item(3, 2).NumberFormat = "@"
$q4Sheet.Cells.Item(3, 2).Value = "180020"
$q4Sheet.Cells.Item(3, 3).Value = "银华成长先锋混合"
$q4Sheet.Cells.Item(3, 4).NumberFormat = "@"
$q4Sheet.Cells.Item(3, 4).Value = "2.18"
$q4Sheet.Cells.Item(3, 5).NumberFormat = "@"
$q4Sheet.Cells.Item(3, 5).Value = "78.61"
$q4Sheet.Cells.Item(3, 6).NumberFormat = "@"
$q4Sheet.Cells.Item(3, 6).Value = "8.46"
$q4Sheet.Cells.Item(3, 7).NumberFormat = "@"
$q4Sheet.Cells.Item(3, 7).Value = "0.1844"
$q4Sheet.Cells.Item(3, 8).Value = 3

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q4 above
#    the existing 2022-Q3 row, then renumber the running index column.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Cells.Item(3, 1).Copy($totalSheet.Cells.Item(2, 1))

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 2
$totalSheet.Cells.Item(2, 4).Value = 1.21

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(7, 1).Value = 5

$excel.CutCopyMode = $false

# Restore the originally-active tab ("总计", the first sheet) now that all
# the edits are done, instead of leaving the newly-inserted sheet selected.
$totalSheet.Activate()
